$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows of data to append (date serial, hora, preco, site, cor)
$data = @(
    @(45215, "19:33", 2100, "amazon", "preto"),
    @(45217, "19:34", 2100, "amazon", "preto"),
    @(45218, "21:26", 2100, "amazon", "preto"),
    @(45220, "10:27", 1954, "amazon", "preto")
)

$startRow = 23
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    $ws.Cells.Item($row, 1).Value = $vals[0]
    $ws.Cells.Item($row, 1).NumberFormat = $ws.Cells.Item(22, 1).NumberFormat

    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}
